$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" column (E) for the four worker rows: bump the period
# from 2508 to 2509, and center the value (matching the rest of the row).
$rng = $ws.Range("E16:E19")
$rng.Value = "2509"
$rng.HorizontalAlignment = -4108  # xlCenter
